$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to force numeric-looking strings to remain text
# (mirrors pasting formatted-as-text content, avoiding Excel's automatic
# numeric conversion when assigning such strings directly to .Value).
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '27.367.89'
$ws.Cells.Item(2, 5).Value = '  +1.24%  '
$ws.Cells.Item(3, 4).Value = '1.824.42'
$ws.Cells.Item(3, 5).Value = '  -0.15%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$helper.Value = '313.67'
$helper.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(5, 5).Value = '  +0.57%  '
$ws.Cells.Item(6, 5).Value = '  +0.01%  '
$helper.Value = '0.4472'
$helper.Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(7, 5).Value = '  +3.29%  '
$helper.Value = '0.3751'
$helper.Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(8, 5).Value = '  +1.98%  '
$helper.Value = '0.07496'
$helper.Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(9, 5).Value = '  +2.90%  '
$helper.Value = '0.8872'
$helper.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(10, 5).Value = '  +4.87%  '
$helper.Value = '21.02'
$helper.Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(11, 5).Value = '  +1.54%  '
$ws.Cells.Item(12, 4).Value = '1.825.48'
$ws.Cells.Item(12, 5).Value = '  -0.10%  '
$helper.Value = '6.757'
$helper.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(13, 5).Value = '  +1.36%  '
$helper.Value = '93.92'
$helper.Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(14, 5).Value = '  +4.91%  '
$helper.Value = '5.407'
$helper.Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(15, 5).Value = '  +2.05%  '
$helper.Value = '0.07105'
$helper.Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(16, 5).Value = '  +0.61%  '
$ws.Cells.Item(17, 5).Value = '  -0.04%  '
$helper.Value = '0.000008805'
$helper.Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(18, 5).Value = '  +0.23%  '
$ws.Cells.Item(19, 5).Value = '  +0.04%  '
$ws.Cells.Item(20, 5).Value = '  +1.99%  '
$ws.Cells.Item(21, 4).Value = '27.379.19'
$ws.Cells.Item(21, 5).Value = '  +0.99%  '
$helper.Value = '5.259'
$helper.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(22, 5).Value = '  +2.13%  '
$helper.Value = '10.91'
$helper.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(23, 5).Value = '  +0.16%  '
$ws.Cells.Item(24, 2).Value = 'Toncoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$helper.Value = '1.963'
$helper.Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(24, 5).Value = '  -1.32%  '
$ws.Cells.Item(25, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$helper.Value = '2.375'
$helper.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(25, 5).Value = '  +7.11%  '
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$helper.Value = '151.47'
$helper.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(26, 5).Value = '  -0.02%  '
$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$helper.Value = '18.62'
$helper.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(27, 5).Value = '  +1.71%  '
$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$helper.Value = '5.356'
$helper.Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(28, 5).Value = '  +2.29%  '
$ws.Cells.Item(29, 2).Value = 'BitcoinCash'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$helper.Value = '118.05'
$helper.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(29, 5).Value = '  +0.78%  '
$ws.Cells.Item(30, 2).Value = 'Stellar'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$helper.Value = '0.08811'
$helper.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(30, 5).Value = '  +0.93%  '
$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$helper.Value = '0.7846'
$helper.Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(31, 5).Value = '  +5.82%  '
$ws.Cells.Item(32, 2).Value = 'ARBITRUM'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$helper.Value = '1.197'
$helper.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(32, 5).Value = '  +1.37%  '
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$helper.Value = '4.514'
$helper.Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(33, 5).Value = '  +1.62%  '
$ws.Cells.Item(34, 2).Value = 'HuobiToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$helper.Value = '2.932'
$helper.Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(34, 5).Value = '  +0.81%  '
$ws.Cells.Item(35, 2).Value = 'Frax'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$helper.Value = '1.001'
$helper.Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(35, 5).Value = '  +0.01%  '
$ws.Cells.Item(36, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$helper.Value = '1.113'
$helper.Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(36, 5).Value = '  +1.38%  '
$ws.Cells.Item(37, 2).Value = 'VeChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$helper.Value = '0.01995'
$helper.Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(37, 5).Value = '  +2.29%  '
$ws.Cells.Item(38, 2).Value = 'Hedera'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$helper.Value = '0.05336'
$helper.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(38, 5).Value = '  +1.69%  '
$ws.Cells.Item(39, 2).Value = 'FraxShare'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$helper.Value = '7.391'
$helper.Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(39, 5).Value = '  +2.30%  '
$ws.Cells.Item(40, 2).Value = 'TheSandbox'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$helper.Value = '0.5316'
$helper.Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(40, 5).Value = '  +3.71%  '
$ws.Cells.Item(41, 2).Value = 'Algorand'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$helper.Value = '0.1727'
$helper.Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(41, 5).Value = '  +1.41%  '
$ws.Cells.Item(42, 2).Value = 'MXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$helper.Value = '2.856'
$helper.Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(42, 5).Value = '  -0.59%  '
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$helper.Value = '2.309'
$helper.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(43, 5).Value = '  +19.27%  '
$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$helper.Value = '8.767'
$helper.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(44, 5).Value = '  +2.04%  '
$ws.Cells.Item(45, 2).Value = 'Decentraland'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$helper.Value = '0.5102'
$helper.Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(45, 5).Value = '  +6.88%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$helper.Value = '10.64'
$helper.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(46, 5).Value = '  +0.12%  '
$ws.Cells.Item(47, 2).Value = 'Quant'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$helper.Value = '105.77'
$helper.Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(47, 5).Value = '  -0.33%  '
$ws.Cells.Item(48, 2).Value = 'NEARProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$helper.Value = '1.701'
$helper.Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(48, 5).Value = '  +2.24%  '
$ws.Cells.Item(49, 2).Value = 'PaxDollar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$helper.Value = '1.001'
$helper.Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(49, 5).Value = '  +0.05%  '
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$helper.Value = '0.06379'
$helper.Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(50, 5).Value = '  +0.62%  '
$ws.Cells.Item(51, 2).Value = 'ThetaToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$helper.Value = '0.9353'
$helper.Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(51, 5).Value = '  +2.81%  '

# Clean up the helper cell so it leaves no trace in the workbook
$helper.Clear() | Out-Null
$excel.CutCopyMode = $false
